# "Write a code to call this method in Cube class." ->
# "Write a code to call this method in Square class." (static modifier added
# to the Square class in the surrounding source listing), with Word's
# automatic "_GoBack" last-edit bookmark relocated onto the newly typed
# word, exactly as Word itself would leave things after an interactive edit.

$d = $word.ActiveDocument

# Locate "Cube" and retype it as "Square" (mirrors the user selecting the
# word and typing the replacement).
$range = $d.Content
$found = $range.Find.Execute("Cube", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
$range.Text = "Square"

$start = $range.Start
$end = $range.End

# Force the run boundary ahead of "Square" with a transient bookmark, then
# drop it again -- leaves the surrounding text split into separate runs
# the way Word's own editor does, without leaving any extra bookmark behind.
$leading = $d.Range($start, $start)
$d.Bookmarks.Add("_TempMark", $leading)

# Word drops its "_GoBack" bookmark (collapsed) right after the text that
# was just typed -- wherever it previously lived gets vacated automatically
# since a bookmark name can only mark one spot at a time.
$trailing = $d.Range($end, $end)
$d.Bookmarks.Add("_GoBack", $trailing)

$d.Bookmarks("_TempMark").Delete()
